$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converters")

# Capture the text of the three comments that live at/after row 3 so we can
# re-create them one row further down (Insert() does not relocate comments).
$commentB3 = $ws.Range("B3").Comment.Text()
$commentB8 = $ws.Range("B8").Comment.Text()
$commentB9 = $ws.Range("B9").Comment.Text()

$ws.Range("B3").Comment.Delete()
$ws.Range("B8").Comment.Delete()
$ws.Range("B9").Comment.Delete()

# Insert a new row above row 3 ("capital cost") for the new
# "fixed capital cost" line item; everything below shifts down one row.
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "fixed capital cost"

# Re-create the moved comments one row down from where they used to be.
$ws.Range("B4").AddComment($commentB3)
$ws.Range("B9").AddComment($commentB8)
$ws.Range("B10").AddComment($commentB9)

# Converters becomes the active sheet/selection, replacing "Time series".
$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
